$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the polyA isolation protocol kit name (affects whole column G, shared string)
$ws.Range("G2:G41").Value = "NEBNextPoly(A)E7490L"

# 2. Re-apply the bold/Arial formatting used by G2:G27 to G28:G41 so the whole
#    column is visually consistent (matches the "accuracy check" cleanup).
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G28:G41").PasteSpecial(-4122) | Out-Null

# 3. Correct the rna sample dates for the last four "2a" rows -- they had all
#    been mistakenly left at "6.5.19"; give each row its own checked date.
foreach ($addr in @("A38", "A39", "A40", "A41")) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("A38").Value = "06.5.19"
$ws.Range("A39").Value = "06.5.20"
$ws.Range("A40").Value = "06.5.21"
$ws.Range("A41").Value = "06.5.22"
foreach ($addr in @("A38", "A39", "A40", "A41")) {
    $ws.Range($addr).NumberFormat = "General"
}

# 4. roboticS1Prep column: turn the literal FALSE values into a live =FALSE()
#    formula for every data row, as part of the accuracy check.
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 9).Formula = "=FALSE()"
}

# 5. Widen the polyAIsolationProtocol column so the longer kit name is visible.
$ws.Columns("G:G").ColumnWidth = 30.5

# 6. Update the view: show gridlines again and leave the selection / scroll
#    position on the newly-corrected rows at the bottom of the sheet.
$excel.ActiveWindow.DisplayGridlines = $true
$ws.Range("A38:A41").Select()
